$wb = $excel.ActiveWorkbook

# "Técnicos" sheet: update the name for técnico id 2 (row 4, column B)
# from "Zé 2" to "Manel 2"
$wsTecnicos = $wb.Worksheets.Item("Técnicos")
$wsTecnicos.Range("B4").Value = "Manel 2"

# Activate the "Técnicos" sheet and select cell B5
$wsTecnicos.Activate()
$wsTecnicos.Range("B5").Select()
